# ------------------------------------------------------------------
# Fixed update to excel issue
#
# 1. Rename header labels on the existing sheets from the generic
#    "Requested quantity" to more specific column names.
# 2. Add a new "PO Forecast" worksheet (after "Monthly Trend") that
#    holds the forecasted PO quantities produced by the forecasting
#    model (ds / PO_Forecast / yhat_lower / yhat_upper).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Rename the "Requested quantity" headers ----------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" worksheet, placed after the last
#     existing sheet ("Monthly Trend") -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# --- Header row -----------------------------------------------------
$ws3.Cells.Item(1,1).Value = "ds"
$ws3.Cells.Item(1,2).Value = "PO_Forecast"
$ws3.Cells.Item(1,3).Value = "yhat_lower"
$ws3.Cells.Item(1,4).Value = "yhat_upper"

# Copy the bold/centered/bordered header style used on the other
# sheets onto the new header row.
$wsWeekly.Range("A1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Cells.Item(2,1).Value = 44934.99999999999
$ws3.Cells.Item(2,2).Value = 27
$ws3.Cells.Item(2,3).Value = -7.859466879760501
$ws3.Cells.Item(2,4).Value = 61.45518748620167
$ws3.Cells.Item(3,1).Value = 44941.99999999999
$ws3.Cells.Item(3,2).Value = 27
$ws3.Cells.Item(3,3).Value = -9.013252782989085
$ws3.Cells.Item(3,4).Value = 59.81274829497256
$ws3.Cells.Item(4,1).Value = 44948.99999999999
$ws3.Cells.Item(4,2).Value = 26
$ws3.Cells.Item(4,3).Value = -9.877070111673635
$ws3.Cells.Item(4,4).Value = 58.3106676224044
$ws3.Cells.Item(5,1).Value = 44955.99999999999
$ws3.Cells.Item(5,2).Value = 26
$ws3.Cells.Item(5,3).Value = -7.608703577376155
$ws3.Cells.Item(5,4).Value = 61.87162481363443
$ws3.Cells.Item(6,1).Value = 44969.99999999999
$ws3.Cells.Item(6,2).Value = 25
$ws3.Cells.Item(6,3).Value = -9.710409855600941
$ws3.Cells.Item(6,4).Value = 58.74481499527687
$ws3.Cells.Item(7,1).Value = 44976.99999999999
$ws3.Cells.Item(7,2).Value = 25
$ws3.Cells.Item(7,3).Value = -6.245608578719798
$ws3.Cells.Item(7,4).Value = 60.2625451315505
$ws3.Cells.Item(8,1).Value = 44990.99999999999
$ws3.Cells.Item(8,2).Value = 24
$ws3.Cells.Item(8,3).Value = -7.587248736824725
$ws3.Cells.Item(8,4).Value = 59.43457334080736
$ws3.Cells.Item(9,1).Value = 44997.99999999999
$ws3.Cells.Item(9,2).Value = 24
$ws3.Cells.Item(9,3).Value = -13.25822567494337
$ws3.Cells.Item(9,4).Value = 55.24832218105306
$ws3.Cells.Item(10,1).Value = 45004.99999999999
$ws3.Cells.Item(10,2).Value = 23
$ws3.Cells.Item(10,3).Value = -10.05146176257137
$ws3.Cells.Item(10,4).Value = 53.53560249274278
$ws3.Cells.Item(11,1).Value = 45102.99999999999
$ws3.Cells.Item(11,2).Value = 18
$ws3.Cells.Item(11,3).Value = -16.26849519156966
$ws3.Cells.Item(11,4).Value = 51.55929724408033
$ws3.Cells.Item(12,1).Value = 45109.99999999999
$ws3.Cells.Item(12,2).Value = 17
$ws3.Cells.Item(12,3).Value = -18.73086466730934
$ws3.Cells.Item(12,4).Value = 51.88106310123747
$ws3.Cells.Item(13,1).Value = 45116.99999999999
$ws3.Cells.Item(13,2).Value = 17
$ws3.Cells.Item(13,3).Value = -18.79976793733841
$ws3.Cells.Item(13,4).Value = 50.11867481124893
$ws3.Cells.Item(14,1).Value = 45137.99999999999
$ws3.Cells.Item(14,2).Value = 16
$ws3.Cells.Item(14,3).Value = -21.51733112947945
$ws3.Cells.Item(14,4).Value = 48.23162656402737
$ws3.Cells.Item(15,1).Value = 45144.99999999999
$ws3.Cells.Item(15,2).Value = 15
$ws3.Cells.Item(15,3).Value = -15.89627367469555
$ws3.Cells.Item(15,4).Value = 48.29351199916555
$ws3.Cells.Item(16,1).Value = 45158.99999999999
$ws3.Cells.Item(16,2).Value = 14
$ws3.Cells.Item(16,3).Value = -20.75579191981328
$ws3.Cells.Item(16,4).Value = 48.58280818949674
$ws3.Cells.Item(17,1).Value = 45165.99999999999
$ws3.Cells.Item(17,2).Value = 14
$ws3.Cells.Item(17,3).Value = -19.90340041045075
$ws3.Cells.Item(17,4).Value = 46.84390462361037
$ws3.Cells.Item(18,1).Value = 45172.99999999999
$ws3.Cells.Item(18,2).Value = 14
$ws3.Cells.Item(18,3).Value = -20.60661712413317
$ws3.Cells.Item(18,4).Value = 49.59659088603564
$ws3.Cells.Item(19,1).Value = 45186.99999999999
$ws3.Cells.Item(19,2).Value = 13
$ws3.Cells.Item(19,3).Value = -19.56113254680215
$ws3.Cells.Item(19,4).Value = 47.49629992835307
$ws3.Cells.Item(20,1).Value = 45193.99999999999
$ws3.Cells.Item(20,2).Value = 13
$ws3.Cells.Item(20,3).Value = -22.01654750982051
$ws3.Cells.Item(20,4).Value = 44.82279798905702
$ws3.Cells.Item(21,1).Value = 45207.99999999999
$ws3.Cells.Item(21,2).Value = 12
$ws3.Cells.Item(21,3).Value = -22.7500355169263
$ws3.Cells.Item(21,4).Value = 44.15650401500737
$ws3.Cells.Item(22,1).Value = 45214.99999999999
$ws3.Cells.Item(22,2).Value = 11
$ws3.Cells.Item(22,3).Value = -21.93087490999667
$ws3.Cells.Item(22,4).Value = 46.21521874229221
$ws3.Cells.Item(23,1).Value = 45221.99999999999
$ws3.Cells.Item(23,2).Value = 11
$ws3.Cells.Item(23,3).Value = -24.42110560384414
$ws3.Cells.Item(23,4).Value = 45.33380417407116
$ws3.Cells.Item(24,1).Value = 45228.99999999999
$ws3.Cells.Item(24,2).Value = 11
$ws3.Cells.Item(24,3).Value = -21.71382147057915
$ws3.Cells.Item(24,4).Value = 45.03969915865344
$ws3.Cells.Item(25,1).Value = 45235.99999999999
$ws3.Cells.Item(25,2).Value = 10
$ws3.Cells.Item(25,3).Value = -23.53141996987424
$ws3.Cells.Item(25,4).Value = 42.53500276819859
$ws3.Cells.Item(26,1).Value = 45242.99999999999
$ws3.Cells.Item(26,2).Value = 10
$ws3.Cells.Item(26,3).Value = -24.31350362776088
$ws3.Cells.Item(26,4).Value = 42.42415893279321
$ws3.Cells.Item(27,1).Value = 45249.99999999999
$ws3.Cells.Item(27,2).Value = 9
$ws3.Cells.Item(27,3).Value = -24.83923331736835
$ws3.Cells.Item(27,4).Value = 41.32349133371185
$ws3.Cells.Item(28,1).Value = 45256.99999999999
$ws3.Cells.Item(28,2).Value = 9
$ws3.Cells.Item(28,3).Value = -27.75041177026052
$ws3.Cells.Item(28,4).Value = 42.13301190845603
$ws3.Cells.Item(29,1).Value = 45263.99999999999
$ws3.Cells.Item(29,2).Value = 9
$ws3.Cells.Item(29,3).Value = -25.57717207581939
$ws3.Cells.Item(29,4).Value = 41.480575037649

# Copy the date-time number-format style used for the "ds" column on
# the other sheets onto the new sheet's A2:A29 range.
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A29").PasteSpecial(-4122)

$excel.CutCopyMode = $false
